$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) values per the latest symbol-list refresh.
$priceUpdates = @{
    "D2" = "266.26"
    "D4" = "6.202"
    "D5" = "0.06146"
    "D6" = "3.560"
    "D7" = "6.715"
    "D8" = "1.361"
    "D9" = "0.8255"
    "D12" = "0.08223"
    "D13" = "0.03400"
    "D14" = "0.03158"
    "D16" = "0.09225"
    "D17" = "3.890"
    "D18" = "0.001694"
    "D19" = "0.04802"
    "D20" = "0.006210"
    "D21" = "0.006289"
    "D22" = "0.001099"
    "D24" = "3.719"
    "D25" = "2.230"
    "D26" = "0.3380"
    "D27" = "0.0002681"
    "D40" = "0.04618"
    "D41" = "0.006986"
    "D42" = "0.1134"
    "D43" = "0.003245"
    "D44" = "0.01103"
    "D45" = "0.00006136"
    "D46" = "0.00000000750"
    "D47" = "0.7700"
    "D48" = "0.2061"
    "D49" = "0.00002101"
    "D50" = "0.01240"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text format so the numeric-looking string is preserved exactly
    # (same representation as the inline string already in the sheet),
    # including trailing zeros, instead of Excel coercing it to a Double.
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cellRef]
}
